$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.000000000000009992007221626409
$ws.Range("B3").Value = -0.000000000000009547918011776346
$ws.Range("B4").Value = -0.000000000000004440892098500626
$ws.Range("B5").Value = 0.000000000000008160139230994901
$ws.Range("B6").Value = 0.000000000000002997602166487923
$ws.Range("B7").Value = 0.00000000000001021405182655144
$ws.Range("B8").Value = -0.00000000000002442490654175344
$ws.Range("B9").Value = 0.00000000000002436939539052219
$ws.Range("B10").Value = -0.00000000000001865174681370263
$ws.Range("B11").Value = -0.000000000000009103828801926284
$ws.Range("B12").Value = 0.00000000000002264854970235319
$ws.Range("B13").Value = 1.593775310996169
$ws.Range("B14").Value = -0.4433754434956397
$ws.Range("B15").Value = -0.5659094841436607
$ws.Range("B16").Value = -0.02440516240001389
$ws.Range("B17").Value = 0.5808185939316554
$ws.Range("B18").Value = 0.08546502771525721
$ws.Range("B19").Value = 0.1413663207244666
$ws.Range("B20").Value = -0.9925327920106826
$ws.Range("B21").Value = 0.6166013166907478
$ws.Range("B22").Value = 0.1324346314480873
$ws.Range("B23").Value = -0.1095386862888253
$ws.Range("B24").Value = -0.2383424898041266
$ws.Range("B25").Value = 0.6794258297128992
$ws.Range("B26").Value = 0.03106364948285656
$ws.Range("B27").Value = -0.4886483869543554
$ws.Range("B28").Value = 0.605758692296735
$ws.Range("B29").Value = 0.2658730860130905
$ws.Range("B30").Value = -0.215129482578987
$ws.Range("B31").Value = -0.4485771038079694
$ws.Range("B32").Value = 0.3419104562980941
$ws.Range("B33").Value = 0.1184663661808556
$ws.Range("B34").Value = 0.4226346167906511
$ws.Range("B35").Value = -0.8394267946364184
$ws.Range("B36").Value = 0.9539990704374277
$ws.Range("B37").Value = -0.1016472839247859
$ws.Range("B38").Value = 0.3894274519612926
$ws.Range("B39").Value = -0.4632625443669942
$ws.Range("B40").Value = -0.005792791899517091
$ws.Range("B41").Value = 0.5624759907254965
$ws.Range("B42").Value = 0.3796720383020056
$ws.Range("B43").Value = -0.3492652477906187
$ws.Range("B44").Value = -0.715889171321582
$ws.Range("B45").Value = 0.1691567084338567
$ws.Range("B46").Value = -0.6525458734825007
$ws.Range("B47").Value = 0.2152667692021368
$ws.Range("B48").Value = -0.885744281990282
$ws.Range("B49").Value = 0.3505440551774829
$ws.Range("B50").Value = 0.4918016298087824
$ws.Range("B51").Value = 0.6329113189306338
$ws.Range("B52").Value = -1.755361674695274
$ws.Range("B53").Value = 0.2210295811551402
$ws.Range("B54").Value = -0.1411005845595663
$ws.Range("B55").Value = 0.05762549224104685
$ws.Range("B56").Value = 0.639181144805368
$ws.Range("B57").Value = -0.1144610931020827
$ws.Range("B58").Value = -0.4184709022183099
$ws.Range("B59").Value = 0.2429549719812211
$ws.Range("B60").Value = 1.128237710175259
$ws.Range("B61").Value = -0.2082144774214872
$ws.Range("B62").Value = -0.6910508056900682
$ws.Range("B63").Value = 0.0947112025147417
$ws.Range("B64").Value = 0.5952133710329237
$ws.Range("B65").Value = -0.1142735084420681
$ws.Range("B66").Value = 0.4225165845712571
$ws.Range("B67").Value = -0.1967006994430545
$ws.Range("B68").Value = 0.8880062631184649
$ws.Range("B69").Value = -0.5897473044240461
$ws.Range("B70").Value = -0.1661013782684851
$ws.Range("B71").Value = -0.3687592045110376
$ws.Range("B72").Value = -0.7537458141733666
$ws.Range("B73").Value = -0.8114147307874704
$ws.Range("B74").Value = -0.1212771608465997
$ws.Range("B75").Value = -0.5935793456314062
$ws.Range("B76").Value = 0.3762390903770712
$ws.Range("B77").Value = -0.4008196643353814
$ws.Range("B78").Value = 0.1721949556597677
$ws.Range("B79").Value = -0.4748495211968715
$ws.Range("B80").Value = 0.9462663862015542
$ws.Range("B81").Value = -0.1356194181296091
$ws.Range("B82").Value = 0.02323282939418525
$ws.Range("B83").Value = -0.5550776869004674
$ws.Range("B84").Value = 0.8865631931960101
$ws.Range("B85").Value = 0.07104253540100011
$ws.Range("B86").Value = 0.1075028472288783
$ws.Range("B87").Value = -0.758469243884506
$ws.Range("B88").Value = -0.5661311224158858
$ws.Range("B89").Value = 0.01540602510147604
$ws.Range("B95").Value = -0.1711104671482269
$ws.Range("B96").Value = 0.2949881563205918
$ws.Range("B97").Value = 0.1633555277145124
$ws.Range("B98").Value = 0.03022793029104015
$ws.Range("B99").Value = -0.2928552586432669
$ws.Range("B100").Value = 0.6214648934017475
$ws.Range("B101").Value = -0.4443674518509457
$ws.Range("B102").Value = -0.1898431833400273
$ws.Range("B103").Value = 0.2351372936899813
$ws.Range("B104").Value = -0.2314726160725444
$ws.Range("B105").Value = -0.1398984511979687
$ws.Range("B106").Value = -0.4898424607665044
$ws.Range("B107").Value = -0.0619887154334321
$ws.Range("B108").Value = 0.4897593076320608
$ws.Range("B109").Value = -0.08025008101715425
$ws.Range("B110").Value = -0.2209309121121734
$ws.Range("B111").Value = -0.2447597831269769
$ws.Range("B112").Value = -0.213315487371266
$ws.Range("B113").Value = 0.188451994675966
$ws.Range("B114").Value = -0.1625432761390304
$ws.Range("B115").Value = -0.03627262647648394
$ws.Range("B116").Value = -0.6870776560732526
$ws.Range("B117").Value = -0.09599133664528381
$ws.Range("B118").Value = 0.02834167208306626
$ws.Range("B119").Value = -0.4067771049315971
$ws.Range("B120").Value = 0.8473145294465119
$ws.Range("B121").Value = -0.3343650428035478
$ws.Range("B122").Value = 0.4558229525073971
$ws.Range("B123").Value = -0.4284561131352164
$ws.Range("B124").Value = -0.2317013822699299
$ws.Range("B125").Value = -0.1321993960636193
$ws.Range("B126").Value = -0.2990191236833212
$ws.Range("B127").Value = 0.3619968567220385
$ws.Range("B128").Value = -0.1919337072516534
$ws.Range("B129").Value = -0.8168372431843682
$ws.Range("B130").Value = -0.1962295261482101
$ws.Range("B131").Value = -0.6265044734712123
$ws.Range("B132").Value = -1.191162284383693
$ws.Range("B133").Value = -0.7174192902096188
$ws.Range("B134").Value = 1.607786042869044
$ws.Range("B135").Value = -0.8848240492139481
$ws.Range("B136").Value = -0.4425935756942813
$ws.Range("B137").Value = -0.01202960564142147
$ws.Range("B138").Value = 0.1065604846180273
$ws.Range("B139").Value = -0.08591989013840001
$ws.Range("B140").Value = 0.1385721760938039
$ws.Range("B141").Value = 0.09472847613588289
$ws.Range("B142").Value = -0.1296176279974082
$ws.Range("B143").Value = -0.01074155887864159
$ws.Range("B144").Value = 0.3186980753357052
$ws.Range("B145").Value = 0.1914876003089772
$ws.Range("B146").Value = 0.08603368373087023
$ws.Range("B147").Value = -0.3828165493744078
$ws.Range("B148").Value = 0.1054308279183608

$ws.Range("A148").Copy()
$ws.Range("A149").PasteSpecial(-4122)
$ws.Range("A149").Value = 45748
$ws.Range("B149").Value = 0

Write-Host "Applied all updates"
